$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G12: mirror F12 value (9->12) with a new number-format style
$ws.Range("G12").Value = $ws.Range("F12").Value

# G14: new value "10->12"
$ws.Range("G14").Value = "10->12"

# G16: new value "Meeting 4", same style as F16/C16/D16
$ws.Range("G16").Value = "Meeting 4"
$ws.Range("G16").Style = $ws.Range("F16").Style

# F17 / G17: long free-text notes, wrapped, left/top aligned
$ws.Range("F17").Value = "Discussed final proposal with Jonty`nWhat is required next week`nStart readding for methadologies report"
$ws.Range("G17").Value = "Discussion about TAS team where everyone was at and what people are doing"

$leftTopWrap = $ws.Range("F17:G17")
$leftTopWrap.HorizontalAlignment = -4131  # xlLeft
$leftTopWrap.VerticalAlignment = -4160    # xlTop
$leftTopWrap.WrapText = $true

# Column widths for F and G
$ws.Columns("F").ColumnWidth = 17.28515625
$ws.Columns("G").ColumnWidth = 15.7109375

# Give G12 the date-ish number format (numFmtId 16 => "d-mmm") to match the new style
$ws.Range("G12").NumberFormat = "d-mmm"

# Update the view: scroll position and selection
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("I17").Select()
